$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
  "10000" = "AKT3"
  "5894"  = "RAF1"
  "2033"  = "EP300"
  "57492" = "ARID1B"
  "6310"  = "ATXN1"
  "2776"  = "GNAQ"
  "23389" = "MED13L"
}

for ($r = 2; $r -le 61; $r++) {
  $cell = $ws.Cells.Item($r, 9)
  $val = $cell.Value2
  if ($val -ne $null) {
    $parts = [string]$val -split "/"
    $newParts = @()
    foreach ($p in $parts) {
      $newParts += $map[$p]
    }
    $cell.Value2 = [string]::Join("/", $newParts)
  }
}
